$d = $word.ActiveDocument

# The "COMPETENCES TECHNIQUES" section contains 7 consecutive single-line
# paragraphs (all sharing identical paragraph formatting) whose order /
# content needs to be updated per the diff. Rather than physically moving
# paragraphs (which carries a risk of disturbing formatting/run
# properties), we rewrite the text content of each of the 7 slots in
# place, since every paragraph here uses the exact same pPr/run
# formatting - this produces output identical to a reorder.

$newTexts = @(
    "Langages : r, python, matlab, c, c++",
    "Visualisation : tableau",
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
    "Autres : dess",
    "Web : client",
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
    "Bases de données : SQL, MongoDB, Neo4j, Redis"
)

# Locate the block: find the paragraph whose text is exactly "Web : client"
# immediately followed by "Langages : ..." etc., by scanning the
# paragraphs collection for the first one whose trimmed content equals
# "Web : client".
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.Replace([char]13, "").Replace([char]7, "")
    if ($t -eq "Web : client") {
        $startIndex = $i
        break
    }
}

for ($j = 0; $j -lt $newTexts.Length; $j++) {
    $p = $d.Paragraphs.Item($startIndex + $j)
    $r = $p.Range
    [void]$r.MoveEnd(1, -1)
    $r.Text = $newTexts[$j]
}
